$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-05-12 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-13 Monday", 2)

$d.Content.Find.Execute("79÷5=15, 4", $true, $false, $false, $false, $false, $true, 1, $false, "14÷2=7, 0", 2)
$d.Content.Find.Execute("50÷5=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "69÷2=34, 1", 2)
$d.Content.Find.Execute("30÷3=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "24÷2=12, 0", 2)
$d.Content.Find.Execute("17÷7=2, 3", $true, $false, $false, $false, $false, $true, 1, $false, "54÷7=7, 5", 2)
$d.Content.Find.Execute("90÷9=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "21÷7=3, 0", 2)

$d.Content.Find.Execute("68÷7=9, 5", $true, $false, $false, $false, $false, $true, 1, $false, "21÷2=10, 1", 2)
$d.Content.Find.Execute("81÷9=9, 0", $true, $false, $false, $false, $false, $true, 1, $false, "13÷2=6, 1", 2)
$d.Content.Find.Execute("57÷2=28, 1", $true, $false, $false, $false, $false, $true, 1, $false, "26÷6=4, 2", 2)
$d.Content.Find.Execute("10÷5=2, 0", $true, $false, $false, $false, $false, $true, 1, $false, "87÷5=17, 2", 2)
$d.Content.Find.Execute("89÷5=17, 4", $true, $false, $false, $false, $false, $true, 1, $false, "41÷3=13, 2", 2)

$d.Content.Find.Execute("31÷2=15, 1", $true, $false, $false, $false, $false, $true, 1, $false, "22÷9=2, 4", 2)
$d.Content.Find.Execute("57÷6=9, 3", $true, $false, $false, $false, $false, $true, 1, $false, "29÷9=3, 2", 2)
$d.Content.Find.Execute("71÷8=8, 7", $true, $false, $false, $false, $false, $true, 1, $false, "12÷2=6, 0", 2)
$d.Content.Find.Execute("58÷2=29, 0", $true, $false, $false, $false, $false, $true, 1, $false, "20÷8=2, 4", 2)
$d.Content.Find.Execute("48÷8=6, 0", $true, $false, $false, $false, $false, $true, 1, $false, "75÷4=18, 3", 2)

$d.Content.Find.Execute("55÷2=27, 1", $true, $false, $false, $false, $false, $true, 1, $false, "93÷7=13, 2", 2)
$d.Content.Find.Execute("79÷6=13, 1", $true, $false, $false, $false, $false, $true, 1, $false, "56÷8=7, 0", 2)
$d.Content.Find.Execute("49÷5=9, 4", $true, $false, $false, $false, $false, $true, 1, $false, "24÷7=3, 3", 2)
$d.Content.Find.Execute("96÷9=10, 6", $true, $false, $false, $false, $false, $true, 1, $false, "62÷8=7, 6", 2)
$d.Content.Find.Execute("79÷9=8, 7", $true, $false, $false, $false, $false, $true, 1, $false, "89÷3=29, 2", 2)

$d.Content.Find.Execute("91÷3=30, 1", $true, $false, $false, $false, $false, $true, 1, $false, "14÷2=7, 0", 2)
$d.Content.Find.Execute("98÷5=19, 3", $true, $false, $false, $false, $false, $true, 1, $false, "56÷3=18, 2", 2)
$d.Content.Find.Execute("76÷5=15, 1", $true, $false, $false, $false, $false, $true, 1, $false, "62÷7=8, 6", 2)
$d.Content.Find.Execute("59÷3=19, 2", $true, $false, $false, $false, $false, $true, 1, $false, "20÷5=4, 0", 2)
$d.Content.Find.Execute("99÷6=16, 3", $true, $false, $false, $false, $false, $true, 1, $false, "72÷6=12, 0", 2)
